# Update the Metadata sheet: bump version, date, and contact info.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# Add a new worksheet "Include from FSIII 2" after the existing "Include from FSIII" sheet,
# mirroring its layout/style but pointing at the new concept (FBOE).
$src = $wb.Worksheets.Item("Include from FSIII")
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "Include from FSIII 2"

$new.Columns.Item(1).ColumnWidth = 30.703125
$new.Columns.Item(2).ColumnWidth = 50.703125

$new.Range("A1").Value = "Property"
$new.Range("B1").Value = "Operation"
$new.Range("C1").Value = "Value"
$new.Range("A1:C1").Style = $src.Range("A1:C1").Style

$new.Range("A2").Value = "concept"
$new.Range("B2").Value = "descendent-of"
$new.Range("C2").Value = "A"
$new.Range("A2:C2").Style = $src.Range("A2:C2").Style

$new.Range("A3").Value = ""
$new.Range("B3").Value = ""
$new.Range("A3:B3").Style = $src.Range("A3:B3").Style

$new.Range("A4").Value = "System URI"
$new.Range("B4").Value = "urn:oid:1.2.208.176.2.21"
$new.Range("A4:B4").Style = $src.Range("A4:B4").Style

# Keep the originally-active "Metadata" tab selected, since adding a sheet makes it active.
$meta.Activate()
